$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the 2023 column (column U) of data, mirroring the style of the
# existing 2022 column (column T) for each row.
$ws.Range("U4").Value = 2023
$ws.Range("U5").Value = 0.5
$ws.Range("U6").Value = 0.3
$ws.Range("U7").Value = 0.4
$ws.Range("U8").Value = 0.4
$ws.Range("U9").Value = 3.2
$ws.Range("U10").Value = 0.6
$ws.Range("U11").Value = "-"
$ws.Range("U12").Value = 0.6
$ws.Range("U13").Value = 0.1
$ws.Range("U14").Value = 0.5

# Copy formatting from the 2022 column (T) into the new 2023 column (U)
# so the new cells match the look of the rest of the table.
$ws.Range("T4:T14").Copy() | Out-Null
$ws.Range("U4:U14").PasteSpecial(-4122) | Out-Null

# Update the selected cell to match the post-edit workbook state.
$ws.Range("B1").Select() | Out-Null
